$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("barnehage")

# New values for columns D (barnehage_antall_plasser) and E (barnehage_ledige_plasser)
# for rows 2-35.
$data = @(
    @{Row=2;  D=38;  E=15},
    @{Row=3;  D=70;  E=67},
    @{Row=4;  D=83;  E=60},
    @{Row=5;  D=61;  E=50},
    @{Row=6;  D=54;  E=50},
    @{Row=7;  D=74;  E=74},
    @{Row=8;  D=92;  E=78},
    @{Row=9;  D=53;  E=53},
    @{Row=10; D=148; E=138},
    @{Row=11; D=98;  E=98},
    @{Row=12; D=47;  E=46},
    @{Row=13; D=73;  E=60},
    @{Row=14; D=60;  E=55},
    @{Row=15; D=59;  E=40},
    @{Row=16; D=66;  E=64},
    @{Row=17; D=35;  E=28},
    @{Row=18; D=82;  E=80},
    @{Row=19; D=96;  E=66},
    @{Row=20; D=51;  E=46},
    @{Row=21; D=40;  E=26},
    @{Row=22; D=206; E=196},
    @{Row=23; D=83;  E=76},
    @{Row=24; D=50;  E=26},
    @{Row=25; D=61;  E=46},
    @{Row=26; D=31;  E=16},
    @{Row=27; D=43;  E=16},
    @{Row=28; D=59;  E=36},
    @{Row=29; D=64;  E=56},
    @{Row=30; D=40;  E=26},
    @{Row=31; D=67;  E=66},
    @{Row=32; D=72;  E=66},
    @{Row=33; D=51;  E=36},
    @{Row=34; D=71;  E=46},
    @{Row=35; D=67;  E=56}
)

foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 4).Value = $entry.D
    $ws.Cells.Item($entry.Row, 5).Value = $entry.E
}
